{"js": "// The source document has the Title, Author and Abstract paragraphs split\n// into many runs (one run per word / space). The edit collapses each of\n// those paragraphs down to a single run holding the full sentence, without\n// touching anything else (table contents, other paragraphs, formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = [\n  { index: 0, text: \"Factsheet: Greek letters\" },\n  { index: 1, text: \"Tom Coleman\" },\n  { index: 3, text: \"Greek letters and their pronunciations in English.\" }\n];\n\nfor (const { index, text } of replacements) {\n  paragraphs.items[index].insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The source document has the Title, Author and Abstract paragraphs split\n# into many runs (one run per word / space). The edit collapses each of\n# those paragraphs down to a single run holding the full sentence, without\n# touching anything else (table contents, other paragraphs, formatting).\n#\n# Word's Find/Replace collapses every run it matches into a single replacement\n# run, which is exactly the \"re-join the split runs\" effect we need, so it is\n# used instead of a plain Range.Text assignment (which only overwrites the\n# first run of a multi-run range).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"Factsheet: Greek letters\" \"Factsheet: Greek letters\"\nReplace-Text \"Tom Coleman\" \"Tom Coleman\"\nReplace-Text \"Greek letters and their pronunciations in English.\" \"Greek letters and their pronunciations in English.\"\n"}
